# Performance.xlsx - "More work JPADCore_v2 concerning the VMC calculation."
#
# Adds a "Minimum control speed (VMC)" row to the TAKE-OFF sheet, just
# below "Stall speed take-off (VsTO)" and above "Decision speed (V1)".
# All the rows that used to follow (Decision speed (V1), Rotation speed
# (V_Rot), Lift-off speed (V_LO), Take-off safety speed (V2), the blank
# separator row and Take-off duration) shift down by one row, keeping
# their original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TAKE-OFF")

# Push rows 10-15 down to 11-16, opening up a blank row 10.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the VMC entry.
$ws.Range("A10").Value = "Minimum control speed (VMC)"
$ws.Range("B10").Value = "m/s"
$ws.Range("C10").Value = 44.07421783106892
